$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet to reflect the new repayment date (2025-09-13)
$ws.Name = "repayment_20250913_20250913"

# --- Row 2 ---
$ws.Cells.Item(2,1).Value = "Annisa Putri Restu"
$ws.Cells.Item(2,2).Value = "Hansyah_S2l"
$ws.Cells.Item(2,3).Value = "S2"
$ws.Cells.Item(2,4).Value = 2
$ws.Cells.Item(2,5).NumberFormat = "@"
$ws.Cells.Item(2,5).Value = "234,046.00"
$ws.Cells.Item(2,5).Style = "Normal"
$ws.Cells.Item(2,6).NumberFormat = "@"
$ws.Cells.Item(2,6).Value = "190,650,342.00"
$ws.Cells.Item(2,6).Style = "Normal"
$ws.Cells.Item(2,7).NumberFormat = "@"
$ws.Cells.Item(2,7).Value = "0.12"
$ws.Cells.Item(2,7).Style = "Normal"
$ws.Cells.Item(2,8).Value = 122
$ws.Cells.Item(2,9).Value = 15
$ws.Cells.Item(2,10).Value = 0
$ws.Cells.Item(2,11).NumberFormat = "@"
$ws.Cells.Item(2,11).Value = "0.00"
$ws.Cells.Item(2,11).Style = "Normal"
$ws.Cells.Item(2,12).NumberFormat = "@"
$ws.Cells.Item(2,12).Value = "0.00"
$ws.Cells.Item(2,12).Style = "Normal"

# --- Row 3 ---
$ws.Cells.Item(3,1).Value = "Yandi Nugraha"
$ws.Cells.Item(3,2).Value = "Hansyah_S2l"
$ws.Cells.Item(3,3).Value = "S2"
$ws.Cells.Item(3,4).Value = 3
$ws.Cells.Item(3,5).NumberFormat = "@"
$ws.Cells.Item(3,5).Value = "1,123,165.00"
$ws.Cells.Item(3,5).Style = "Normal"
$ws.Cells.Item(3,6).NumberFormat = "@"
$ws.Cells.Item(3,6).Value = "151,034,134.00"
$ws.Cells.Item(3,6).Style = "Normal"
$ws.Cells.Item(3,7).NumberFormat = "@"
$ws.Cells.Item(3,7).Value = "0.74"
$ws.Cells.Item(3,7).Style = "Normal"
$ws.Cells.Item(3,8).Value = 7
$ws.Cells.Item(3,9).Value = 15
$ws.Cells.Item(3,10).Value = 0
$ws.Cells.Item(3,11).NumberFormat = "@"
$ws.Cells.Item(3,11).Value = "0.00"
$ws.Cells.Item(3,11).Style = "Normal"
$ws.Cells.Item(3,12).NumberFormat = "@"
$ws.Cells.Item(3,12).Value = "0.00"
$ws.Cells.Item(3,12).Style = "Normal"

# --- Row 4 ---
$ws.Cells.Item(4,1).Value = "Wasti Feronika Sihombing"
$ws.Cells.Item(4,2).Value = "Hansyah_S2l"
$ws.Cells.Item(4,3).Value = "S2"
$ws.Cells.Item(4,4).Value = 2
$ws.Cells.Item(4,5).NumberFormat = "@"
$ws.Cells.Item(4,5).Value = "1,249,987.00"
$ws.Cells.Item(4,5).Style = "Normal"
$ws.Cells.Item(4,6).NumberFormat = "@"
$ws.Cells.Item(4,6).Value = "158,163,068.00"
$ws.Cells.Item(4,6).Style = "Normal"
$ws.Cells.Item(4,7).NumberFormat = "@"
$ws.Cells.Item(4,7).Value = "0.79"
$ws.Cells.Item(4,7).Style = "Normal"
$ws.Cells.Item(4,8).Value = 50
$ws.Cells.Item(4,9).Value = 15
$ws.Cells.Item(4,10).Value = 0
$ws.Cells.Item(4,11).NumberFormat = "@"
$ws.Cells.Item(4,11).Value = "0.00"
$ws.Cells.Item(4,11).Style = "Normal"
$ws.Cells.Item(4,12).NumberFormat = "@"
$ws.Cells.Item(4,12).Value = "0.00"
$ws.Cells.Item(4,12).Style = "Normal"

# --- Row 5 ---
$ws.Cells.Item(5,1).Value = "Axl Wicaksono"
$ws.Cells.Item(5,2).Value = "Hansyah_S2l"
$ws.Cells.Item(5,3).Value = "S2"
$ws.Cells.Item(5,4).Value = 2
$ws.Cells.Item(5,5).NumberFormat = "@"
$ws.Cells.Item(5,5).Value = "706,349.00"
$ws.Cells.Item(5,5).Style = "Normal"
$ws.Cells.Item(5,6).NumberFormat = "@"
$ws.Cells.Item(5,6).Value = "144,104,643.00"
$ws.Cells.Item(5,6).Style = "Normal"
$ws.Cells.Item(5,7).NumberFormat = "@"
$ws.Cells.Item(5,7).Value = "0.49"
$ws.Cells.Item(5,7).Style = "Normal"
$ws.Cells.Item(5,8).Value = 23
$ws.Cells.Item(5,9).Value = 15
$ws.Cells.Item(5,10).Value = 1
$ws.Cells.Item(5,11).NumberFormat = "@"
$ws.Cells.Item(5,11).Value = "3.14"
$ws.Cells.Item(5,11).Style = "Normal"
$ws.Cells.Item(5,12).NumberFormat = "@"
$ws.Cells.Item(5,12).Value = "6.67"
$ws.Cells.Item(5,12).Style = "Normal"

# --- Row 6 ---
$ws.Cells.Item(6,1).Value = "Nuraini"
$ws.Cells.Item(6,2).Value = "Hansyah_S2l"
$ws.Cells.Item(6,3).Value = "S2"
$ws.Cells.Item(6,4).Value = 1
$ws.Cells.Item(6,5).NumberFormat = "@"
$ws.Cells.Item(6,5).Value = "678,767.00"
$ws.Cells.Item(6,5).Style = "Normal"
$ws.Cells.Item(6,6).NumberFormat = "@"
$ws.Cells.Item(6,6).Value = "99,841,865.00"
$ws.Cells.Item(6,6).Style = "Normal"
$ws.Cells.Item(6,7).NumberFormat = "@"
$ws.Cells.Item(6,7).Value = "0.68"
$ws.Cells.Item(6,7).Style = "Normal"
$ws.Cells.Item(6,8).Value = 123
$ws.Cells.Item(6,9).Value = 9
$ws.Cells.Item(6,10).Value = 0
$ws.Cells.Item(6,11).NumberFormat = "@"
$ws.Cells.Item(6,11).Value = "0.00"
$ws.Cells.Item(6,11).Style = "Normal"
$ws.Cells.Item(6,12).NumberFormat = "@"
$ws.Cells.Item(6,12).Value = "0.00"
$ws.Cells.Item(6,12).Style = "Normal"

# --- Row 7 ---
$ws.Cells.Item(7,1).Value = "Riska Nurlita"
$ws.Cells.Item(7,2).Value = "Hansyah_S2l"
$ws.Cells.Item(7,3).Value = "S2"
$ws.Cells.Item(7,4).Value = 1
$ws.Cells.Item(7,5).NumberFormat = "@"
$ws.Cells.Item(7,5).Value = "200,000.00"
$ws.Cells.Item(7,5).Style = "Normal"
$ws.Cells.Item(7,6).NumberFormat = "@"
$ws.Cells.Item(7,6).Value = "192,661,552.00"
$ws.Cells.Item(7,6).Style = "Normal"
$ws.Cells.Item(7,7).NumberFormat = "@"
$ws.Cells.Item(7,7).Value = "0.10"
$ws.Cells.Item(7,7).Style = "Normal"
$ws.Cells.Item(7,8).Value = 238
$ws.Cells.Item(7,9).Value = 16
$ws.Cells.Item(7,10).Value = 0
$ws.Cells.Item(7,11).NumberFormat = "@"
$ws.Cells.Item(7,11).Value = "0.00"
$ws.Cells.Item(7,11).Style = "Normal"
$ws.Cells.Item(7,12).NumberFormat = "@"
$ws.Cells.Item(7,12).Value = "0.00"
$ws.Cells.Item(7,12).Style = "Normal"

# --- Row 8 ---
$ws.Cells.Item(8,1).Value = "Debora Retima Sihombing"
$ws.Cells.Item(8,2).Value = "Hansyah_S2l"
$ws.Cells.Item(8,3).Value = "S2"
$ws.Cells.Item(8,4).Value = 0
$ws.Cells.Item(8,5).NumberFormat = "@"
$ws.Cells.Item(8,5).Value = "0.00"
$ws.Cells.Item(8,5).Style = "Normal"
$ws.Cells.Item(8,6).NumberFormat = "@"
$ws.Cells.Item(8,6).Value = "156,270,040.00"
$ws.Cells.Item(8,6).Style = "Normal"
$ws.Cells.Item(8,7).NumberFormat = "@"
$ws.Cells.Item(8,7).Value = "0.00"
$ws.Cells.Item(8,7).Style = "Normal"
$ws.Cells.Item(8,8).Value = 92
$ws.Cells.Item(8,9).Value = 15
$ws.Cells.Item(8,10).Value = 0
$ws.Cells.Item(8,11).NumberFormat = "@"
$ws.Cells.Item(8,11).Value = "0.00"
$ws.Cells.Item(8,11).Style = "Normal"
$ws.Cells.Item(8,12).NumberFormat = "@"
$ws.Cells.Item(8,12).Value = "0.00"
$ws.Cells.Item(8,12).Style = "Normal"

# --- Row 9 ---
$ws.Cells.Item(9,1).Value = "Azizah Rahmawati"
$ws.Cells.Item(9,2).Value = "Hansyah_S2l"
$ws.Cells.Item(9,3).Value = "S2"
$ws.Cells.Item(9,4).Value = 0
$ws.Cells.Item(9,5).NumberFormat = "@"
$ws.Cells.Item(9,5).Value = "0.00"
$ws.Cells.Item(9,5).Style = "Normal"
$ws.Cells.Item(9,6).NumberFormat = "@"
$ws.Cells.Item(9,6).Value = "172,276,267.00"
$ws.Cells.Item(9,6).Style = "Normal"
$ws.Cells.Item(9,7).NumberFormat = "@"
$ws.Cells.Item(9,7).Value = "0.00"
$ws.Cells.Item(9,7).Style = "Normal"
$ws.Cells.Item(9,8).Value = 62
$ws.Cells.Item(9,9).Value = 15
$ws.Cells.Item(9,10).Value = 0
$ws.Cells.Item(9,11).NumberFormat = "@"
$ws.Cells.Item(9,11).Value = "0.00"
$ws.Cells.Item(9,11).Style = "Normal"
$ws.Cells.Item(9,12).NumberFormat = "@"
$ws.Cells.Item(9,12).Value = "0.00"
$ws.Cells.Item(9,12).Style = "Normal"

# --- Row 10 ---
$ws.Cells.Item(10,1).Value = "Erlangga Hutama"
$ws.Cells.Item(10,2).Value = "Hansyah_S2l"
$ws.Cells.Item(10,3).Value = "S2"
$ws.Cells.Item(10,4).Value = 0
$ws.Cells.Item(10,5).NumberFormat = "@"
$ws.Cells.Item(10,5).Value = "0.00"
$ws.Cells.Item(10,5).Style = "Normal"
$ws.Cells.Item(10,6).NumberFormat = "@"
$ws.Cells.Item(10,6).Value = "135,862,450.00"
$ws.Cells.Item(10,6).Style = "Normal"
$ws.Cells.Item(10,7).NumberFormat = "@"
$ws.Cells.Item(10,7).Value = "0.00"
$ws.Cells.Item(10,7).Style = "Normal"
$ws.Cells.Item(10,8).Value = 0
$ws.Cells.Item(10,9).Value = 15
$ws.Cells.Item(10,10).Value = 0
$ws.Cells.Item(10,11).NumberFormat = "@"
$ws.Cells.Item(10,11).Value = "0.00"
$ws.Cells.Item(10,11).Style = "Normal"
$ws.Cells.Item(10,12).NumberFormat = "@"
$ws.Cells.Item(10,12).Value = "0.00"
$ws.Cells.Item(10,12).Style = "Normal"

# --- Row 11 ---
$ws.Cells.Item(11,1).Value = "Erick Ervan Dewanggga"
$ws.Cells.Item(11,2).Value = "Hansyah_S2l"
$ws.Cells.Item(11,3).Value = "S2"
$ws.Cells.Item(11,4).Value = 0
$ws.Cells.Item(11,5).NumberFormat = "@"
$ws.Cells.Item(11,5).Value = "0.00"
$ws.Cells.Item(11,5).Style = "Normal"
$ws.Cells.Item(11,6).NumberFormat = "@"
$ws.Cells.Item(11,6).Value = "151,560,437.00"
$ws.Cells.Item(11,6).Style = "Normal"
$ws.Cells.Item(11,7).NumberFormat = "@"
$ws.Cells.Item(11,7).Value = "0.00"
$ws.Cells.Item(11,7).Style = "Normal"
$ws.Cells.Item(11,8).Value = 128
$ws.Cells.Item(11,9).Value = 15
$ws.Cells.Item(11,10).Value = 0
$ws.Cells.Item(11,11).NumberFormat = "@"
$ws.Cells.Item(11,11).Value = "0.00"
$ws.Cells.Item(11,11).Style = "Normal"
$ws.Cells.Item(11,12).NumberFormat = "@"
$ws.Cells.Item(11,12).Value = "0.00"
$ws.Cells.Item(11,12).Style = "Normal"

# --- Row 12 ---
$ws.Cells.Item(12,1).Value = "Ridhoi Berkat Zebua"
$ws.Cells.Item(12,2).Value = "Hansyah_S2l"
$ws.Cells.Item(12,3).Value = "S2"
$ws.Cells.Item(12,4).Value = 0
$ws.Cells.Item(12,5).NumberFormat = "@"
$ws.Cells.Item(12,5).Value = "0.00"
$ws.Cells.Item(12,5).Style = "Normal"
$ws.Cells.Item(12,6).NumberFormat = "@"
$ws.Cells.Item(12,6).Value = "162,481,134.00"
$ws.Cells.Item(12,6).Style = "Normal"
$ws.Cells.Item(12,7).NumberFormat = "@"
$ws.Cells.Item(12,7).Value = "0.00"
$ws.Cells.Item(12,7).Style = "Normal"
$ws.Cells.Item(12,8).Value = 85
$ws.Cells.Item(12,9).Value = 15
$ws.Cells.Item(12,10).Value = 0
$ws.Cells.Item(12,11).NumberFormat = "@"
$ws.Cells.Item(12,11).Value = "0.00"
$ws.Cells.Item(12,11).Style = "Normal"
$ws.Cells.Item(12,12).NumberFormat = "@"
$ws.Cells.Item(12,12).Value = "0.00"
$ws.Cells.Item(12,12).Style = "Normal"

# --- Row 13 ---
$ws.Cells.Item(13,1).Value = "Romli"
$ws.Cells.Item(13,2).Value = "Hansyah_S2l"
$ws.Cells.Item(13,3).Value = "S2"
$ws.Cells.Item(13,4).Value = 0
$ws.Cells.Item(13,5).NumberFormat = "@"
$ws.Cells.Item(13,5).Value = "0.00"
$ws.Cells.Item(13,5).Style = "Normal"
$ws.Cells.Item(13,6).NumberFormat = "@"
$ws.Cells.Item(13,6).Value = "147,555,539.00"
$ws.Cells.Item(13,6).Style = "Normal"
$ws.Cells.Item(13,7).NumberFormat = "@"
$ws.Cells.Item(13,7).Value = "0.00"
$ws.Cells.Item(13,7).Style = "Normal"
$ws.Cells.Item(13,8).Value = 10
$ws.Cells.Item(13,9).Value = 15
$ws.Cells.Item(13,10).Value = 0
$ws.Cells.Item(13,11).NumberFormat = "@"
$ws.Cells.Item(13,11).Value = "0.00"
$ws.Cells.Item(13,11).Style = "Normal"
$ws.Cells.Item(13,12).NumberFormat = "@"
$ws.Cells.Item(13,12).Value = "0.00"
$ws.Cells.Item(13,12).Style = "Normal"

# --- Row 14 ---
$ws.Cells.Item(14,1).Value = "Fadilah Damayanti"
$ws.Cells.Item(14,2).Value = "Hansyah_S2l"
$ws.Cells.Item(14,3).Value = "S2"
$ws.Cells.Item(14,4).Value = 0
$ws.Cells.Item(14,5).NumberFormat = "@"
$ws.Cells.Item(14,5).Value = "0.00"
$ws.Cells.Item(14,5).Style = "Normal"
$ws.Cells.Item(14,6).NumberFormat = "@"
$ws.Cells.Item(14,6).Value = "149,841,389.00"
$ws.Cells.Item(14,6).Style = "Normal"
$ws.Cells.Item(14,7).NumberFormat = "@"
$ws.Cells.Item(14,7).Value = "0.00"
$ws.Cells.Item(14,7).Style = "Normal"
$ws.Cells.Item(14,8).Value = 0
$ws.Cells.Item(14,9).Value = 15
$ws.Cells.Item(14,10).Value = 0
$ws.Cells.Item(14,11).NumberFormat = "@"
$ws.Cells.Item(14,11).Value = "0.00"
$ws.Cells.Item(14,11).Style = "Normal"
$ws.Cells.Item(14,12).NumberFormat = "@"
$ws.Cells.Item(14,12).Value = "0.00"
$ws.Cells.Item(14,12).Style = "Normal"

# --- Row 15 ---
$ws.Cells.Item(15,1).Value = "Aldi Taufik"
$ws.Cells.Item(15,2).Value = "Hansyah_S2l"
$ws.Cells.Item(15,3).Value = "S2"
$ws.Cells.Item(15,4).Value = 2
$ws.Cells.Item(15,5).NumberFormat = "@"
$ws.Cells.Item(15,5).Value = "195,000.00"
$ws.Cells.Item(15,5).Style = "Normal"
$ws.Cells.Item(15,6).NumberFormat = "@"
$ws.Cells.Item(15,6).Value = "168,949,279.00"
$ws.Cells.Item(15,6).Style = "Normal"
$ws.Cells.Item(15,7).NumberFormat = "@"
$ws.Cells.Item(15,7).Value = "0.12"
$ws.Cells.Item(15,7).Style = "Normal"
$ws.Cells.Item(15,8).Value = 108
$ws.Cells.Item(15,9).Value = 15
$ws.Cells.Item(15,10).Value = 0
$ws.Cells.Item(15,11).NumberFormat = "@"
$ws.Cells.Item(15,11).Value = "0.00"
$ws.Cells.Item(15,11).Style = "Normal"
$ws.Cells.Item(15,12).NumberFormat = "@"
$ws.Cells.Item(15,12).Value = "0.00"
$ws.Cells.Item(15,12).Style = "Normal"

# --- Row 16 ---
$ws.Cells.Item(16,1).Value = "Nur Halim"
$ws.Cells.Item(16,2).Value = "Hansyah_S2l"
$ws.Cells.Item(16,3).Value = "S2"
$ws.Cells.Item(16,4).Value = 0
$ws.Cells.Item(16,5).NumberFormat = "@"
$ws.Cells.Item(16,5).Value = "0.00"
$ws.Cells.Item(16,5).Style = "Normal"
$ws.Cells.Item(16,6).NumberFormat = "@"
$ws.Cells.Item(16,6).Value = "150,970,103.00"
$ws.Cells.Item(16,6).Style = "Normal"
$ws.Cells.Item(16,7).NumberFormat = "@"
$ws.Cells.Item(16,7).Value = "0.00"
$ws.Cells.Item(16,7).Style = "Normal"
$ws.Cells.Item(16,8).Value = 281
$ws.Cells.Item(16,9).Value = 15
$ws.Cells.Item(16,10).Value = 0
$ws.Cells.Item(16,11).NumberFormat = "@"
$ws.Cells.Item(16,11).Value = "0.00"
$ws.Cells.Item(16,11).Style = "Normal"
$ws.Cells.Item(16,12).NumberFormat = "@"
$ws.Cells.Item(16,12).Value = "0.00"
$ws.Cells.Item(16,12).Style = "Normal"

# --- Row 17 ---
$ws.Cells.Item(17,1).Value = "Adistira Winditya P"
$ws.Cells.Item(17,2).Value = "Hansyah_S2l"
$ws.Cells.Item(17,3).Value = "S2"
$ws.Cells.Item(17,4).Value = 0
$ws.Cells.Item(17,5).NumberFormat = "@"
$ws.Cells.Item(17,5).Value = "0.00"
$ws.Cells.Item(17,5).Style = "Normal"
$ws.Cells.Item(17,6).NumberFormat = "@"
$ws.Cells.Item(17,6).Value = "152,248,724.00"
$ws.Cells.Item(17,6).Style = "Normal"
$ws.Cells.Item(17,7).NumberFormat = "@"
$ws.Cells.Item(17,7).Value = "0.00"
$ws.Cells.Item(17,7).Style = "Normal"
$ws.Cells.Item(17,8).Value = 27
$ws.Cells.Item(17,9).Value = 15
$ws.Cells.Item(17,10).Value = 0
$ws.Cells.Item(17,11).NumberFormat = "@"
$ws.Cells.Item(17,11).Value = "0.00"
$ws.Cells.Item(17,11).Style = "Normal"
$ws.Cells.Item(17,12).NumberFormat = "@"
$ws.Cells.Item(17,12).Value = "0.00"
$ws.Cells.Item(17,12).Style = "Normal"

# --- Row 18 ---
$ws.Cells.Item(18,1).Value = "Sucika Wardani"
$ws.Cells.Item(18,2).Value = "Hansyah_S2l"
$ws.Cells.Item(18,3).Value = "S2"
$ws.Cells.Item(18,4).Value = 7
$ws.Cells.Item(18,5).NumberFormat = "@"
$ws.Cells.Item(18,5).Value = "1,018,178.00"
$ws.Cells.Item(18,5).Style = "Normal"
$ws.Cells.Item(18,6).NumberFormat = "@"
$ws.Cells.Item(18,6).Value = "146,306,782.00"
$ws.Cells.Item(18,6).Style = "Normal"
$ws.Cells.Item(18,7).NumberFormat = "@"
$ws.Cells.Item(18,7).Value = "0.70"
$ws.Cells.Item(18,7).Style = "Normal"
$ws.Cells.Item(18,8).Value = 122
$ws.Cells.Item(18,9).Value = 15
$ws.Cells.Item(18,10).Value = 0
$ws.Cells.Item(18,11).NumberFormat = "@"
$ws.Cells.Item(18,11).Value = "0.00"
$ws.Cells.Item(18,11).Style = "Normal"
$ws.Cells.Item(18,12).NumberFormat = "@"
$ws.Cells.Item(18,12).Value = "0.00"
$ws.Cells.Item(18,12).Style = "Normal"

# Update the active selection to match the latest save state
$ws.Range("A2:A18").Select() | Out-Null
